$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update computed/forecast values for rows 2-7 (columns B, C, E, F)
# Row 2
$ws.Range("B2").Value = 10346.2550516824
$ws.Range("C2").Value = 10486.158009027
$ws.Range("E2").Value = 7678.55671612344
$ws.Range("F2").Value = -42.1477197853997

# Row 3
$ws.Range("B3").Value = 10607.7771522341
$ws.Range("C3").Value = 10764.4780174285
$ws.Range("E3").Value = 8302.95136991572
$ws.Range("F3").Value = 326.465391139344

# Row 4
$ws.Range("B4").Value = 11098.8116464094
$ws.Range("C4").Value = 11554.7146371413
$ws.Range("E4").Value = 8697.91670040976
$ws.Range("F4").Value = 375.848805731293

# Row 5
$ws.Range("C5").Value = 9279.21861270876
$ws.Range("E5").Value = 8884.37040927386
$ws.Range("F5").Value = 253.655375915942

# Row 6
$ws.Range("C6").Value = 11060.0363488547
$ws.Range("E6").Value = 9062.5408570272
$ws.Range("F6").Value = 335.279883578412

# Row 7
$ws.Range("C7").Value = 7809.84776872701
$ws.Range("E7").Value = 8513.25258611682
$ws.Range("F7").Value = 176.968348118493
